$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.359.33"
$ws.Range("E2").Value = "  -1.00%  "

$ws.Range("D3").Value = "3.321.12"
$ws.Range("E3").Value = "  +1.78%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "185.95"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.39%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "577.91"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.83%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.605"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.90%  "

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("E9").Value = "  -0.04%  "

$ws.Range("E10").Value = "  +0.69%  "

$ws.Range("E11").Value = "  -0.13%  "

$ws.Range("D12").Value = "3.898.53"
$ws.Range("E12").Value = "  +1.76%  "

$ws.Range("E13").Value = "  -0.75%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.38"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.47%  "

$ws.Range("D15").Value = "67.553.30"
$ws.Range("E15").Value = "  -0.73%  "

$ws.Range("E16").Value = "  -0.17%  "

$ws.Range("D17").Value = "3.316.05"
$ws.Range("E17").Value = "  +1.23%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "443.37"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +6.60%  "

$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.57"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.11%  "

$ws.Range("B20").Value = "Polkadot"
$ws.Range("C20").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.67"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.58%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.72"
$ws.Range("D21").Style = "Normal"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "74.19"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.02%  "

$ws.Range("E23").Value = "  -0.14%  "

$ws.Range("D24").Value = "3.465.36"
$ws.Range("E24").Value = "  +1.73%  "

$ws.Range("E25").Value = "  +1.18%  "

$ws.Range("E26").Value = "  +1.65%  "

$ws.Range("E27").Value = "  +1.48%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.05"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.87%  "

$ws.Range("E29").Value = "  -0.90%  "

$ws.Range("E30").Value = "  +1.54%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.92"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.17%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.34"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.25%  "

$ws.Range("E33").Value = "  -0.07%  "

$ws.Range("E34").Value = "  -0.30%  "

$ws.Range("E35").Value = "  -0.20%  "

$ws.Range("E36").Value = "  +4.99%  "

$ws.Range("E37").Value = "  -0.20%  "

$ws.Range("B38").Value = "EnergySwap"
$ws.Range("C38").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "27.27"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.96%  "

$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.84"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.11%  "

$ws.Range("D40").Value = "2.788.27"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.789"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.86%  "

$ws.Range("E42").Value = "  +0.60%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.23"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.81%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "40.32"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.13%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0671"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.59%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "24.70"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.46%  "

$ws.Range("E47").Value = "  -1.29%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "327.26"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.89%  "

$ws.Range("E49").Value = "  -0.27%  "

$ws.Range("E50").Value = "  +1.37%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "31.12"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.66%  "
